# Auto-generated edit script applying numeric/text updates
# from the cosinor-per-day CircaDiPy re-run (square_10, period=7).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = [double]"25.57000000000056"
$ws.Range("H2").Value = [double]"4.740652315149418e-13"
$ws.Range("I2").Value = [double]"4.740652315149418e-13"
$ws.Range("L2").Value = [double]"56.73207308611578"
$ws.Range("M2").Value = "[44.913499294097335, 68.55064687813423]"
$ws.Range("N2").Value = [double]"1.480593425640109e-12"
$ws.Range("O2").Value = [double]"1.480593425640109e-12"
$ws.Range("P2").Value = [double]"1.616395018964118"
$ws.Range("Q2").Value = "[1.3773949772495797, 1.8553950606786564]"
$ws.Range("T2").Value = [double]"50.8009259535134"
$ws.Range("U2").Value = "[42.97859883663241, 58.62325307039439]"
$ws.Range("X2").Value = [double]"18.99193193193235"
$ws.Range("Y2").Value = [double]"18.01929929929969"
$ws.Range("Z2").Value = [double]"19.964564564565"
# Row 3
$ws.Range("F3").Value = [double]"25.57000000000056"
$ws.Range("H3").Value = [double]"2.834399381868025e-13"
$ws.Range("I3").Value = [double]"2.834399381868025e-13"
$ws.Range("L3").Value = [double]"58.22600740904101"
$ws.Range("M3").Value = "[44.576535937076684, 71.87547888100534]"
$ws.Range("N3").Value = [double]"4.821631982565577e-11"
$ws.Range("O3").Value = [double]"4.821631982565577e-11"
$ws.Range("P3").Value = [double]"1.17613178422681"
$ws.Range("Q3").Value = "[0.9245527929483472, 1.4277107755052718]"
$ws.Range("R3").Value = [double]"3.309352791802667e-12"
$ws.Range("S3").Value = [double]"3.309352791802667e-12"
$ws.Range("T3").Value = [double]"55.18732805422803"
$ws.Range("U3").Value = "[47.45842946911834, 62.91622663933771]"
$ws.Range("X3").Value = [double]"20.78362362362408"
$ws.Range("Y3").Value = [double]"19.75979979980023"
$ws.Range("Z3").Value = [double]"21.80744744744792"
# Row 4
$ws.Range("F4").Value = [double]"25.57000000000056"
$ws.Range("H4").Value = [double]"6.364908600176022e-13"
$ws.Range("I4").Value = [double]"6.364908600176022e-13"
$ws.Range("L4").Value = [double]"59.69277641234851"
$ws.Range("M4").Value = "[44.022061125566125, 75.36349169913089]"
$ws.Range("N4").Value = [double]"1.038726882285346e-09"
$ws.Range("O4").Value = [double]"1.038726882285346e-09"
$ws.Range("P4").Value = [double]"0.748447499053424"
$ws.Range("Q4").Value = "[0.4842895582110396, 1.0126054398958084]"
$ws.Range("R4").Value = [double]"8.494301764194745e-07"
$ws.Range("S4").Value = [double]"8.494301764194745e-07"
$ws.Range("T4").Value = [double]"54.37291424167596"
$ws.Range("U4").Value = "[46.471896371023774, 62.273932112328154]"
$ws.Range("X4").Value = [double]"22.52412412412462"
$ws.Range("Y4").Value = [double]"21.44910910910958"
$ws.Range("Z4").Value = [double]"23.59913913913965"
# Row 5
$ws.Range("F5").Value = [double]"25.57000000000056"
$ws.Range("H5").Value = [double]"2.564615186884112e-14"
$ws.Range("I5").Value = [double]"2.564615186884112e-14"
$ws.Range("L5").Value = [double]"59.89025986033145"
$ws.Range("M5").Value = "[46.953364186493275, 72.82715553416962]"
$ws.Range("N5").Value = [double]"4.443334589154802e-12"
$ws.Range("O5").Value = [double]"4.443334589154802e-12"
$ws.Range("P5").Value = [double]"0.3836579616996545"
$ws.Range("Q5").Value = "[0.15723686954903915, 0.6100790538502698]"
$ws.Range("R5").Value = [double]"0.001369965683592156"
$ws.Range("S5").Value = [double]"0.001369965683592156"
$ws.Range("T5").Value = [double]"54.01044940187423"
$ws.Range("U5").Value = "[46.92643660509089, 61.09446219865756]"
$ws.Range("X5").Value = [double]"24.00866866866919"
$ws.Range("Y5").Value = [double]"23.08722722722773"
$ws.Range("Z5").Value = [double]"24.93011011011065"
# Row 6
$ws.Range("F6").Value = [double]"25.57000000000056"
$ws.Range("H6").Value = [double]"6.559086607182962e-12"
$ws.Range("I6").Value = [double]"6.559086607182962e-12"
$ws.Range("L6").Value = [double]"59.59209256395634"
$ws.Range("M6").Value = "[46.429229848659546, 72.75495527925314]"
$ws.Range("N6").Value = [double]"8.628209258176867e-12"
$ws.Range("O6").Value = [double]"8.628209258176867e-12"
$ws.Range("P6").Value = [double]"-0.02515789912784605"
$ws.Range("Q6").Value = "[-0.2641579408423853, 0.2138421425866932]"
$ws.Range("R6").Value = [double]"0.8330559143396632"
$ws.Range("S6").Value = [double]"0.8330559143396632"
$ws.Range("T6").Value = [double]"55.4091153424128"
$ws.Range("U6").Value = "[46.997854939194795, 63.8203757456308]"
$ws.Range("V6").Value = [double]"0"
$ws.Range("W6").Value = [double]"0"
$ws.Range("X6").Value = [double]"0.1023823823823875"
$ws.Range("Y6").Value = [double]"-0.870250250250268"
$ws.Range("Z6").Value = [double]"1.075015015015043"
# Row 7
$ws.Range("F7").Value = [double]"25.57000000000056"
$ws.Range("H7").Value = [double]"4.951816734433123e-12"
$ws.Range("I7").Value = [double]"4.951816734433123e-12"
$ws.Range("L7").Value = [double]"55.27475014780693"
$ws.Range("M7").Value = "[40.628120290875415, 69.92138000473844]"
$ws.Range("N7").Value = [double]"1.320848985741918e-09"
$ws.Range("O7").Value = [double]"1.320848985741918e-09"
$ws.Range("P7").Value = [double]"-0.3522105877898465"
$ws.Range("Q7").Value = "[-0.6289474781961548, -0.07547369738353815]"
$ws.Range("R7").Value = [double]"0.01377767788085027"
$ws.Range("S7").Value = [double]"0.01377767788085027"
$ws.Range("T7").Value = [double]"57.15570471464797"
$ws.Range("U7").Value = "[49.29158406362688, 65.01982536566905]"
$ws.Range("V7").Value = [double]"0"
$ws.Range("W7").Value = [double]"0"
$ws.Range("X7").Value = [double]"1.433353353353386"
$ws.Range("Y7").Value = [double]"0.3071471471471543"
$ws.Range("Z7").Value = [double]"2.559559559559617"
# Row 8
$ws.Range("F8").Value = [double]"25.48000000000054"
$ws.Range("H8").Value = [double]"2.02171612784241e-13"
$ws.Range("I8").Value = [double]"2.02171612784241e-13"
$ws.Range("L8").Value = [double]"58.82362411701318"
$ws.Range("M8").Value = "[43.83085041096976, 73.8163978230566]"
$ws.Range("N8").Value = [double]"4.7846904216442e-10"
$ws.Range("O8").Value = [double]"4.7846904216442e-10"
$ws.Range("P8").Value = [double]"-0.9811580659860013"
$ws.Range("Q8").Value = "[-1.2327370572644627, -0.7295790747075399]"
$ws.Range("R8").Value = [double]"5.608367104059653e-10"
$ws.Range("S8").Value = [double]"5.608367104059653e-10"
$ws.Range("T8").Value = [double]"52.64452922070537"
$ws.Range("U8").Value = "[44.80419152411643, 60.48486691729431]"
$ws.Range("V8").Value = [double]"0"
$ws.Range("W8").Value = [double]"0"
$ws.Range("X8").Value = [double]"3.978858858858942"
$ws.Range("Y8").Value = [double]"2.958638638638702"
$ws.Range("Z8").Value = [double]"4.999079079079182"
# Row 9
$ws.Range("F9").Value = [double]"25.48000000000054"
$ws.Range("H9").Value = [double]"9.248157795127554e-14"
$ws.Range("I9").Value = [double]"9.248157795127554e-14"
$ws.Range("L9").Value = [double]"56.67983021903125"
$ws.Range("M9").Value = "[44.01317297234738, 69.34648746571511]"
$ws.Range("N9").Value = [double]"1.216404754700307e-11"
$ws.Range("O9").Value = [double]"1.216404754700307e-11"
$ws.Range("P9").Value = [double]"-1.283052855520156"
$ws.Range("Q9").Value = "[-1.534631846798618, -1.0314738642416934]"
$ws.Range("R9").Value = [double]"2.229327833447314e-13"
$ws.Range("S9").Value = [double]"2.229327833447314e-13"
$ws.Range("T9").Value = [double]"52.74898427560244"
$ws.Range("U9").Value = "[45.31866707223075, 60.179301478974125]"
$ws.Range("V9").Value = [double]"0"
$ws.Range("W9").Value = [double]"0"
$ws.Range("X9").Value = [double]"5.203123123123234"
$ws.Range("Y9").Value = [double]"4.18290290290299"
$ws.Range("Z9").Value = [double]"6.223343343343478"
# Row 10
$ws.Range("F10").Value = [double]"25.48000000000054"
$ws.Range("H10").Value = [double]"2.65631960871815e-12"
$ws.Range("I10").Value = [double]"2.65631960871815e-12"
$ws.Range("L10").Value = [double]"54.06985840827898"
$ws.Range("M10").Value = "[42.456267023272815, 65.68344979328513]"
$ws.Range("N10").Value = [double]"3.747668841924678e-12"
$ws.Range("O10").Value = [double]"3.747668841924678e-12"
$ws.Range("P10").Value = [double]"-1.559789745926464"
$ws.Range("Q10").Value = "[-1.8239476867688493, -1.2956318050840787]"
$ws.Range("R10").Value = [double]"1.77635683940025e-15"
$ws.Range("S10").Value = [double]"1.77635683940025e-15"
$ws.Range("T10").Value = [double]"59.59637290098173"
$ws.Range("U10").Value = "[51.75204750363848, 67.44069829832497]"
$ws.Range("X10").Value = [double]"6.325365365365503"
$ws.Range("Y10").Value = [double]"5.254134134134247"
$ws.Range("Z10").Value = [double]"7.396596596596758"
# Row 11
$ws.Range("F11").Value = [double]"25.48000000000054"
$ws.Range("H11").Value = [double]"1.012745443063068e-12"
$ws.Range("I11").Value = [double]"1.012745443063068e-12"
$ws.Range("L11").Value = [double]"56.28273788870311"
$ws.Range("M11").Value = "[42.64185058307257, 69.92362519433365]"
$ws.Range("N11").Value = [double]"1.224029766433432e-10"
$ws.Range("O11").Value = [double]"1.224029766433432e-10"
$ws.Range("P11").Value = [double]"-1.899421384152387"
$ws.Range("Q11").Value = "[-2.176158274558696, -1.6226844937460791]"
$ws.Range("R11").Value = [double]"0"
$ws.Range("S11").Value = [double]"0"
$ws.Range("T11").Value = [double]"53.32896268759693"
$ws.Range("U11").Value = "[45.526862356807655, 61.1310630183862]"
$ws.Range("X11").Value = [double]"7.702662662662828"
$ws.Range("Y11").Value = [double]"6.580420420420561"
$ws.Range("Z11").Value = [double]"8.824904904905095"
# Row 12
$ws.Range("F12").Value = [double]"25.48000000000054"
$ws.Range("H12").Value = [double]"1.936895088761048e-11"
$ws.Range("I12").Value = [double]"1.936895088761048e-11"
$ws.Range("L12").Value = [double]"59.62722081300951"
$ws.Range("M12").Value = "[42.116011759959704, 77.13842986605931]"
$ws.Range("N12").Value = [double]"1.6541970015993e-08"
$ws.Range("O12").Value = [double]"1.6541970015993e-08"
$ws.Range("P12").Value = [double]"-2.377421467581465"
$ws.Range("Q12").Value = "[-2.6667373075516965, -2.0881056276112337]"
$ws.Range("T12").Value = [double]"51.91963294684471"
$ws.Range("U12").Value = "[43.08997945475725, 60.74928643893218]"
$ws.Range("V12").Value = [double]"1.998401444325282e-15"
$ws.Range("W12").Value = [double]"1.998401444325282e-15"
$ws.Range("X12").Value = [double]"9.641081081081287"
$ws.Range("Y12").Value = [double]"8.467827827828009"
$ws.Range("Z12").Value = [double]"10.81433433433457"
# Row 13
$ws.Range("F13").Value = [double]"25.48000000000054"
$ws.Range("H13").Value = [double]"4.907185768843192e-14"
$ws.Range("I13").Value = [double]"4.907185768843192e-14"
$ws.Range("L13").Value = [double]"61.13888268616628"
$ws.Range("M13").Value = "[47.40987034144669, 74.86789503088586]"
$ws.Range("N13").Value = [double]"1.4000578474338e-11"
$ws.Range("O13").Value = [double]"1.4000578474338e-11"
$ws.Range("P13").Value = [double]"-2.779947853627005"
$ws.Range("Q13").Value = "[-3.006368945777621, -2.5535267614763884]"
$ws.Range("T13").Value = [double]"51.16816375026361"
$ws.Range("U13").Value = "[43.7820491318091, 58.554278368718116]"
$ws.Range("X13").Value = [double]"11.27343343343368"
$ws.Range("Y13").Value = [double]"10.35523523523545"
$ws.Range("Z13").Value = [double]"12.1916316316319"
# Row 14
$ws.Range("F14").Value = [double]"25.48000000000054"
$ws.Range("H14").Value = [double]"1.339233168806686e-10"
$ws.Range("I14").Value = [double]"1.339233168806686e-10"
$ws.Range("L14").Value = [double]"59.35201016771131"
$ws.Range("M14").Value = "[45.087593016463515, 73.6164273189591]"
$ws.Range("N14").Value = [double]"9.698597480678472e-11"
$ws.Range("O14").Value = [double]"9.698597480678472e-11"
$ws.Range("P14").Value = [double]"-3.119579491852928"
$ws.Range("Q14").Value = "[-3.3963163822592364, -2.8428426014466197]"
$ws.Range("T14").Value = [double]"53.69464956094838"
$ws.Range("U14").Value = "[44.43375637551122, 62.95554274638555]"
$ws.Range("V14").Value = [double]"3.108624468950438e-15"
$ws.Range("W14").Value = [double]"3.108624468950438e-15"
$ws.Range("X14").Value = [double]"12.650730730731"
$ws.Range("Y14").Value = [double]"11.52848848848874"
$ws.Range("Z14").Value = [double]"13.77297297297327"
